# Mice_table.xlsx edit
#
# The lab corrected its Genotype vocabulary: the generic labels
# "Heterozygous" / "Wildtype" are replaced by the specific allele
# names "R403Q(+/-)" and "Null(-) " throughout the Genotype column.
#
# Cells that previously carried the (red-highlighted) "Heterozygous"
# label were retyped by hand and lost their red-font formatting, so
# after the edit every such cell uses the plain/default cell style
# (the same style already used in D2/D3). The "Wildtype" cells kept
# whatever formatting they already had.
#
# Also bring the selection/scroll position to where the edit was last
# made (cell D7 on row 7, with the view scrolled down a bit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole used range of the sheet.
$allData = $ws.Range("A1:G53")

# --- Replace the old genotype labels with the new allele names -------------
# Using whole-cell matching (xlWhole = 1) so we only touch cells whose
# entire content is exactly "Heterozygous" / "Wildtype" (the shared-string
# table has no other cell containing these words as a substring).
$allData.Replace("Heterozygous", "R403Q(+/-)", 1) | Out-Null
$allData.Replace("Wildtype", "Null(-) ", 1) | Out-Null

# --- Restore plain styling on the retyped "Heterozygous" -> "R403Q(+/-)" ---
# cells that used to carry the red highlight style (the ones in the
# "G#" mouse rows use style index 8 / red font); copy the plain style
# already used by D2 (a cell that was always in the default style) onto
# them so they match the rest of the column instead of keeping the old
# red highlight.
$ws.Range("D2").Copy() | Out-Null
$plainStyleTargets = @(14, 15, 17, 18, 21, 22, 23, 24, 25)
foreach ($row in $plainStyleTargets) {
    $ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = $false

# --- Move the on-screen selection / scroll position -------------------------
$win = $excel.ActiveWindow
try { $win.ScrollRow = 3 } catch {}
try { $win.ScrollColumn = 1 } catch {}
$ws.Range("D7").Select()
